# Progress.xlsx edit: mark row "1c" (row 11) as fully done.
# - D11: percentage progress 0 -> 100
# - E11: "Encarregue" (assignee) "-" -> "Bernardo"
# - F11: status formula recalculates TODO -> Done! automatically
# - I3: average formula recalculates automatically
# - Selection moves to J13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the progress percentage for row 11 ("1c")
$ws.Range("D11").Value = 100

# Assign "Bernardo" as responsible for row 11
$ws.Range("E11").Value = "Bernardo"

# Move the active selection to J13, as recorded in the sheet view
$ws.Range("J13").Select()

$wb.Save()
